$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Cells.Item(1, 4).Value = "MontoBruto"
$ws.Cells.Item(1, 5).Value = "MedioPago"
$ws.Cells.Item(1, 6).Value = "DiagnósticoBot"

# --- New data table values ---
$ids = @(18293000, 18293001, 18293002, 18293003, 18293004, 18293005, 18293006, 18293007, 18293008)
$fecha = "02/12/2025"
$productos = @("Gaseosa 500ml", "Galletitas", "Caramelo", "Chocolate", "Chicle", "Agua mineral 500ml", "Snack salado", "Cigarrillos", "Alfajor")
$montos = @(3750.0, 2080.0, 200.0, 890.0, 320.0, 3000.0, 2800.0, 12500.0, 650.0)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]

    # "02/12/2025" parses as a real date, so force text first or Excel will
    # silently store a date serial instead of the literal string.
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $fecha
    $cellB.ClearFormats()

    $ws.Cells.Item($row, 3).Value = $productos[$i]
    $ws.Cells.Item($row, 4).Value = $montos[$i]
}

# --- Remove old extra rows (11 and 12) ---
$ws.Range("A11:F12").Clear()
